$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 939.7917
$ws.Range("I2").Value = 155.2
$ws.Range("J2").Value = 1146.2632
$ws.Range("K2").Value = 155.2
$ws.Range("L2").Value = 1146.2632
$ws.Range("M2").Value = -42.19999999999999
$ws.Range("N2").Value = -1372.2632
$ws.Range("H6").Value = 1159.75
$ws.Range("I6").Value = 1318.2
$ws.Range("K6").Value = 3954.6
$ws.Range("M6").Value = -3842.6
$ws.Range("H11").Value = 31
$ws.Range("I11").Value = 31
$ws.Range("K11").Value = 31
$ws.Range("M11").Value = 109
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H17").Value = 2639
$ws.Range("J17").Value = 2639
$ws.Range("L17").Value = 7917
$ws.Range("N17").Value = -8253
$ws.Range("H32").Value = 3654.9092
$ws.Range("I32").Value = 3851.6667
$ws.Range("J32").Value = 3581.125
$ws.Range("K32").Value = 3851.6667
$ws.Range("L32").Value = 3581.125
$ws.Range("M32").Value = -3525.6667
$ws.Range("N32").Value = -4233.125
$ws.Range("H39").Value = 1506.2858
$ws.Range("I39").Value = 90.666664
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 271.999992
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = 24.00000799999998
$ws.Range("N39").Value = -30592
$ws.Range("H40").Value = 15083565
$ws.Range("I40").Value = 7143817
$ws.Range("J40").Value = 25008250
$ws.Range("K40").Value = 7143817
$ws.Range("L40").Value = 25008250
$ws.Range("M40").Value = -7143642
$ws.Range("N40").Value = -25008600
$ws.Range("H41").Value = 498.2143
$ws.Range("I41").Value = 468.8889
$ws.Range("J41").Value = 551
$ws.Range("K41").Value = 468.8889
$ws.Range("L41").Value = 551
$ws.Range("M41").Value = -28.88889999999998
$ws.Range("N41").Value = -1431
$ws.Range("H45").Value = 12785.714
$ws.Range("I45").Value = 8333.333000000001
$ws.Range("K45").Value = 24999.999
$ws.Range("M45").Value = -24807.999
$ws.Range("H52").Value = 5332.6665
$ws.Range("I52").Value = 2999
$ws.Range("J52").Value = 10000
$ws.Range("K52").Value = 8997
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -8837
$ws.Range("N52").Value = -30320
$ws.Range("H58").Value = 20462.889
$ws.Range("I58").Value = 115
$ws.Range("J58").Value = 26276.572
$ws.Range("K58").Value = 345
$ws.Range("L58").Value = 78829.716
$ws.Range("M58").Value = -195
$ws.Range("N58").Value = -79129.716
$ws.Range("H59").Value = 1000
$ws.Range("I59").Value = 1000
$ws.Range("K59").Value = 3000
$ws.Range("M59").Value = -2443
$ws.Range("H61").Value = 1553.75
$ws.Range("I61").Value = 1553.75
$ws.Range("K61").Value = 4661.25
$ws.Range("M61").Value = -4489.25
$ws.Range("H62").Value = 7567.4287
$ws.Range("I62").Value = 8095.8
$ws.Range("J62").Value = 6246.5
$ws.Range("K62").Value = 8095.8
$ws.Range("L62").Value = 6246.5
$ws.Range("M62").Value = -7471.8
$ws.Range("N62").Value = -7494.5
$ws.Range("H64").Value = 6833.1113
$ws.Range("I64").Value = 6098.2856
$ws.Range("J64").Value = 7300.727
$ws.Range("K64").Value = 6098.2856
$ws.Range("L64").Value = 7300.727
$ws.Range("M64").Value = -5850.2856
$ws.Range("N64").Value = -7796.727
$ws.Range("H65").Value = 7567.4287
$ws.Range("I65").Value = 8095.8
$ws.Range("J65").Value = 6246.5
$ws.Range("K65").Value = 40479
$ws.Range("L65").Value = 31232.5
$ws.Range("M65").Value = -37359
$ws.Range("N65").Value = -37472.5
$ws.Range("H67").Value = 6833.1113
$ws.Range("I67").Value = 6098.2856
$ws.Range("J67").Value = 7300.727
$ws.Range("K67").Value = 6098.2856
$ws.Range("L67").Value = 7300.727
$ws.Range("M67").Value = -5240.2856
$ws.Range("N67").Value = -9016.726999999999
$ws.Range("H69").Value = 7444.385
$ws.Range("J69").Value = 8480
$ws.Range("L69").Value = 25440
$ws.Range("N69").Value = -27188
$ws.Range("H70").Value = 4777.5557
$ws.Range("I70").Value = 5100
$ws.Range("J70").Value = 4132.6665
$ws.Range("K70").Value = 15300
$ws.Range("L70").Value = 12397.9995
$ws.Range("M70").Value = -15030
$ws.Range("N70").Value = -12937.9995
$ws.Range("H72").Value = 7444.385
$ws.Range("J72").Value = 8480
$ws.Range("L72").Value = 76320
$ws.Range("N72").Value = -85056
$ws.Range("H73").Value = 4777.5557
$ws.Range("I73").Value = 5100
$ws.Range("J73").Value = 4132.6665
$ws.Range("K73").Value = 15300
$ws.Range("L73").Value = 12397.9995
$ws.Range("M73").Value = -14364
$ws.Range("N73").Value = -14269.9995
$ws.Range("H96").Value = 595.75
$ws.Range("I96").Value = 627.2727
$ws.Range("J96").Value = 249
$ws.Range("K96").Value = 1881.8181
$ws.Range("L96").Value = 747
$ws.Range("M96").Value = -508.8181
$ws.Range("N96").Value = -3493
$ws.Range("H98").Value = 2536
$ws.Range("I98").Value = 1012.8333
$ws.Range("K98").Value = 1012.8333
$ws.Range("M98").Value = 485.1667
$ws.Range("H106").Value = 12121
$ws.Range("I106").Value = 12308.75
$ws.Range("K106").Value = 12308.75
$ws.Range("M106").Value = -11677.75
$ws.Range("H116").Value = 3371.2646
$ws.Range("I116").Value = 3136
$ws.Range("J116").Value = 3469.2917
$ws.Range("K116").Value = 3136
$ws.Range("L116").Value = 3469.2917
$ws.Range("M116").Value = 306
$ws.Range("N116").Value = -10353.2917
$ws.Range("H122").Value = 2536
$ws.Range("I122").Value = 1012.8333
$ws.Range("K122").Value = 3038.4999
$ws.Range("M122").Value = -588.4998999999998
$ws.Range("H137").Value = 41679456
$ws.Range("I137").Value = 62517510
$ws.Range("J137").Value = 3336.75
$ws.Range("K137").Value = 187552530
$ws.Range("L137").Value = 10010.25
$ws.Range("M137").Value = -187549980
$ws.Range("N137").Value = -15110.25
$ws.Range("H138").Value = 5333.615
$ws.Range("I138").Value = 4442.4443
$ws.Range("J138").Value = 5600.967
$ws.Range("K138").Value = 13327.3329
$ws.Range("L138").Value = 16802.901
$ws.Range("M138").Value = -8187.332900000001
$ws.Range("N138").Value = -27082.901

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 250
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -134
$ws.Range("N4").ClearContents()
$ws.Range("H14").Value = 6950
$ws.Range("I14").Value = 2500
$ws.Range("J14").Value = 8433.333000000001
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 8433.333000000001
$ws.Range("M14").Value = -2325
$ws.Range("N14").Value = -8783.333000000001
$ws.Range("H32").Value = 12457.25
$ws.Range("I32").Value = 8441.385
$ws.Range("J32").Value = 22898.5
$ws.Range("K32").Value = 8441.385
$ws.Range("L32").Value = 22898.5
$ws.Range("M32").Value = -8154.385
$ws.Range("N32").Value = -23472.5
$ws.Range("H37").Value = 500020030
$ws.Range("I37").Value = 500020030
$ws.Range("K37").Value = 500020030
$ws.Range("M37").Value = -500019757
$ws.Range("H74").Value = 18896.166
$ws.Range("I74").Value = 17792.666
$ws.Range("K74").Value = 17792.666
$ws.Range("M74").Value = -16918.666
$ws.Range("H77").Value = 18896.166
$ws.Range("I77").Value = 17792.666
$ws.Range("K77").Value = 88963.33
$ws.Range("M77").Value = -84595.33
$ws.Range("H97").Value = 610.8182
$ws.Range("J97").Value = 968.1667
$ws.Range("L97").Value = 968.1667
$ws.Range("N97").Value = -1960.1667
$ws.Range("H110").Value = 2016
$ws.Range("I110").Value = 1613.5625
$ws.Range("K110").Value = 1613.5625
$ws.Range("M110").Value = 431.4375
$ws.Range("H132").Value = 9451.352000000001
$ws.Range("I132").Value = 6864.4614
$ws.Range("K132").Value = 20593.3842
$ws.Range("M132").Value = -18063.3842
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2948
$ws.Range("I20").Value = 2640.8667
$ws.Range("K20").Value = 2640.8667
$ws.Range("M20").Value = -2393.8667
$ws.Range("H76").Value = 19749.75
$ws.Range("J76").Value = 19749.75
$ws.Range("L76").Value = 19749.75
$ws.Range("N76").Value = -20379.75
$ws.Range("H79").Value = 19749.75
$ws.Range("J79").Value = 19749.75
$ws.Range("L79").Value = 19749.75
$ws.Range("N79").Value = -21933.75
$ws.Range("H105").Value = 1366.9231
$ws.Range("I105").Value = 1337.619
$ws.Range("J105").Value = 1490
$ws.Range("K105").Value = 1337.619
$ws.Range("L105").Value = 1490
$ws.Range("M105").Value = 409.3810000000001
$ws.Range("N105").Value = -4984
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H134").Value = 3314.2
$ws.Range("I134").Value = 1190.3334
$ws.Range("K134").Value = 3571.0002
$ws.Range("M134").Value = -1036.0002
$ws.Range("H141").Value = 249998.5
$ws.Range("J141").Value = 249998.5
$ws.Range("L141").Value = 249998.5
$ws.Range("N141").Value = -260358.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 528231.3
$ws.Range("I31").Value = 11451.583
$ws.Range("J31").Value = 941655.0600000001
$ws.Range("K31").Value = 11451.583
$ws.Range("L31").Value = 941655.0600000001
$ws.Range("M31").Value = -11156.583
$ws.Range("N31").Value = -942245.0600000001
$ws.Range("H34").Value = 528231.3
$ws.Range("I34").Value = 11451.583
$ws.Range("J34").Value = 941655.0600000001
$ws.Range("K34").Value = 11451.583
$ws.Range("L34").Value = 941655.0600000001
$ws.Range("M34").Value = -11249.583
$ws.Range("N34").Value = -942059.0600000001
$ws.Range("H105").Value = 10147.678
$ws.Range("I105").Value = 11282.777
$ws.Range("K105").Value = 11282.777
$ws.Range("M105").Value = -9535.777
$ws.Range("H122").Value = 2036.4546
$ws.Range("I122").Value = 1885.421
$ws.Range("J122").Value = 2993
$ws.Range("K122").Value = 5656.263
$ws.Range("L122").Value = 8979
$ws.Range("M122").Value = -3206.263
$ws.Range("N122").Value = -13879
$ws.Range("H132").Value = 2640.3125
$ws.Range("I132").Value = 2088.926
$ws.Range("J132").Value = 5617.8
$ws.Range("K132").Value = 6266.778
$ws.Range("L132").Value = 16853.4
$ws.Range("M132").Value = -3736.778
$ws.Range("N132").Value = -21913.4
$ws.Range("H134").Value = 1689.6666
$ws.Range("I134").Value = 1850
$ws.Range("J134").Value = 1506.4286
$ws.Range("K134").Value = 5550
$ws.Range("L134").Value = 4519.2858
$ws.Range("M134").Value = -3015
$ws.Range("N134").Value = -9589.2858

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 200129.62
$ws.Range("J4").Value = 60067
$ws.Range("K4").Value = 600388.86
$ws.Range("L4").Value = 180201
$ws.Range("M4").Value = -600276.86
$ws.Range("N4").Value = -180425
$ws.Range("H5").Value = 694.875
$ws.Range("I5").Value = 694.875
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2084.625
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1972.625
$ws.Range("N5").ClearContents()
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H8").Value = 583.6429000000001
$ws.Range("I8").Value = 583.6429000000001
$ws.Range("K8").Value = 1750.9287
$ws.Range("M8").Value = -1611.9287
$ws.Range("H9").Value = 350
$ws.Range("I9").Value = 350
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1050
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -826
$ws.Range("N9").ClearContents()
$ws.Range("H10").Value = 133.66667
$ws.Range("I10").Value = 100.5
$ws.Range("K10").Value = 301.5
$ws.Range("M10").Value = -162.5
$ws.Range("H11").Value = 2351
$ws.Range("I11").Value = 570
$ws.Range("J11").Value = 3114.2856
$ws.Range("K11").Value = 1710
$ws.Range("L11").Value = 9342.856800000001
$ws.Range("M11").Value = -1570
$ws.Range("N11").Value = -9622.856800000001
$ws.Range("H29").Value = 293.5
$ws.Range("I29").Value = 297
$ws.Range("J29").Value = 290
$ws.Range("K29").Value = 891
$ws.Range("L29").Value = 870
$ws.Range("M29").Value = -614
$ws.Range("N29").Value = -1424
$ws.Range("H39").Value = 6613.4
$ws.Range("J39").Value = 6356.6665
$ws.Range("L39").Value = 19069.9995
$ws.Range("N39").Value = -19657.9995
$ws.Range("H47").Value = 1181.6666
$ws.Range("I47").Value = 22.5
$ws.Range("K47").Value = 67.5
$ws.Range("M47").Value = 363.5
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H80").Value = 10000
$ws.Range("I80").Value = 10000
$ws.Range("K80").Value = 30000
$ws.Range("M80").Value = -29064
$ws.Range("H83").Value = 10000
$ws.Range("I83").Value = 10000
$ws.Range("K83").Value = 90000
$ws.Range("M83").Value = -85320
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 2500
$ws.Range("K99").Value = 7500
$ws.Range("M99").Value = -5254
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H104").Value = 10000
$ws.Range("J104").Value = 10000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -35242
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H114").Value = 18596.416
$ws.Range("I114").Value = 376.75
$ws.Range("J114").Value = 27706.25
$ws.Range("K114").Value = 1130.25
$ws.Range("L114").Value = 83118.75
$ws.Range("M114").Value = 2123.75
$ws.Range("N114").Value = -89626.75
$ws.Range("H120").Value = 15800
$ws.Range("J120").Value = 15800
$ws.Range("L120").Value = 47400
$ws.Range("N120").Value = -57076
$ws.Range("H121").Value = 167832.17
$ws.Range("J121").Value = 1331.6666
$ws.Range("L121").Value = 3994.9998
$ws.Range("N121").Value = -6614.9998
$ws.Range("H129").Value = 7058.1665
$ws.Range("I129").Value = 8375.857
$ws.Range("J129").Value = 6219.636
$ws.Range("K129").Value = 25127.571
$ws.Range("L129").Value = 18658.908
$ws.Range("M129").Value = -20127.571
$ws.Range("N129").Value = -28658.908
$ws.Range("H131").Value = 2727.7812
$ws.Range("I131").Value = 1583.3636
$ws.Range("J131").Value = 3327.238
$ws.Range("K131").Value = 4750.0908
$ws.Range("L131").Value = 9981.714
$ws.Range("M131").Value = 289.9092000000001
$ws.Range("N131").Value = -20061.714
$ws.Range("H132").Value = 4338.4
$ws.Range("I132").Value = 1849.5
$ws.Range("J132").Value = 5997.6665
$ws.Range("K132").Value = 16645.5
$ws.Range("L132").Value = 53978.9985
$ws.Range("M132").Value = -14115.5
$ws.Range("N132").Value = -59038.9985
$ws.Range("H135").Value = 694.875
$ws.Range("I135").Value = 694.875
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6253.875
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3718.875
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 5899.4
$ws.Range("J137").Value = 6249.25
$ws.Range("L137").Value = 18747.75
$ws.Range("N137").Value = -28947.75
$ws.Range("H139").Value = 24139.348
$ws.Range("I139").Value = 25678.2
$ws.Range("K139").Value = 77034.60000000001
$ws.Range("M139").Value = -71894.60000000001
$ws.Range("H140").Value = 2059.5
$ws.Range("I140").Value = 1448.6923
$ws.Range("K140").Value = 4346.0769
$ws.Range("M140").Value = 833.9231

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 6451
$ws.Range("J55").Value = 6451
$ws.Range("L55").Value = 6451
$ws.Range("N55").Value = -7105
$ws.Range("H58").Value = 9527.625
$ws.Range("I58").Value = 5555.25
$ws.Range("J58").Value = 13500
$ws.Range("K58").Value = 5555.25
$ws.Range("L58").Value = 13500
$ws.Range("M58").Value = -5278.25
$ws.Range("N58").Value = -14054
$ws.Range("H70").Value = 18617.863
$ws.Range("I70").Value = 26048.23
$ws.Range("J70").Value = 7885.1113
$ws.Range("K70").Value = 26048.23
$ws.Range("L70").Value = 7885.1113
$ws.Range("M70").Value = -25778.23
$ws.Range("N70").Value = -8425.1113
$ws.Range("H73").Value = 18617.863
$ws.Range("I73").Value = 26048.23
$ws.Range("J73").Value = 7885.1113
$ws.Range("K73").Value = 26048.23
$ws.Range("L73").Value = 7885.1113
$ws.Range("M73").Value = -25112.23
$ws.Range("N73").Value = -9757.1113
$ws.Range("H80").Value = 2749.5
$ws.Range("I80").Value = 2332.6667
$ws.Range("K80").Value = 2332.6667
$ws.Range("M80").Value = -1334.6667
$ws.Range("H83").Value = 2749.5
$ws.Range("I83").Value = 2332.6667
$ws.Range("K83").Value = 11663.3335
$ws.Range("M83").Value = -6671.333500000001
$ws.Range("H122").Value = 4753.8125
$ws.Range("I122").Value = 3417.4
$ws.Range("J122").Value = 6981.1665
$ws.Range("K122").Value = 10252.2
$ws.Range("L122").Value = 20943.4995
$ws.Range("M122").Value = -7802.200000000001
$ws.Range("N122").Value = -25843.4995
$ws.Range("H126").Value = 4515
$ws.Range("I126").Value = 4744.1055
$ws.Range("J126").Value = 4180.154
$ws.Range("K126").Value = 14232.3165
$ws.Range("L126").Value = 12540.462
$ws.Range("M126").Value = -11762.3165
$ws.Range("N126").Value = -17480.462
$ws.Range("H132").Value = 7130.1665
$ws.Range("I132").Value = 7146.45
$ws.Range("J132").Value = 7048.75
$ws.Range("K132").Value = 21439.35
$ws.Range("L132").Value = 21146.25
$ws.Range("M132").Value = -18909.35
$ws.Range("N132").Value = -26206.25
$ws.Range("H139").Value = 149999
$ws.Range("J139").Value = 149999
$ws.Range("L139").Value = 149999
$ws.Range("N139").Value = -160279

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 71434570
$ws.Range("I7").Value = 90914296
$ws.Range("K7").Value = 90914296
$ws.Range("M7").Value = -90914184
$ws.Range("H40").Value = 4641.5835
$ws.Range("I40").Value = 4641.5835
$ws.Range("K40").Value = 4641.5835
$ws.Range("M40").Value = -4505.5835
$ws.Range("H55").Value = 3299.5
$ws.Range("I55").Value = 999.3333
$ws.Range("K55").Value = 999.3333
$ws.Range("M55").Value = -826.3333
$ws.Range("H61").Value = 2644
$ws.Range("I61").Value = 2710.9167
$ws.Range("J61").Value = 1038
$ws.Range("K61").Value = 2710.9167
$ws.Range("L61").Value = 1038
$ws.Range("M61").Value = -2508.9167
$ws.Range("N61").Value = -1442
$ws.Range("H68").Value = 4733.091
$ws.Range("I68").Value = 2612.375
$ws.Range("J68").Value = 10388.333
$ws.Range("K68").Value = 2612.375
$ws.Range("L68").Value = 10388.333
$ws.Range("M68").Value = -1863.375
$ws.Range("N68").Value = -11886.333
$ws.Range("H71").Value = 4733.091
$ws.Range("I71").Value = 2612.375
$ws.Range("J71").Value = 10388.333
$ws.Range("K71").Value = 13061.875
$ws.Range("L71").Value = 51941.665
$ws.Range("M71").Value = -9317.875
$ws.Range("N71").Value = -59429.665
$ws.Range("H82").Value = 4417.5
$ws.Range("I82").Value = 2243
$ws.Range("J82").Value = 7461.8
$ws.Range("K82").Value = 2243
$ws.Range("L82").Value = 7461.8
$ws.Range("M82").Value = -1882
$ws.Range("N82").Value = -8183.8
$ws.Range("H85").Value = 4417.5
$ws.Range("I85").Value = 2243
$ws.Range("J85").Value = 7461.8
$ws.Range("K85").Value = 2243
$ws.Range("L85").Value = 7461.8
$ws.Range("M85").Value = -995
$ws.Range("N85").Value = -9957.799999999999
$ws.Range("H93").Value = 2727.8635
$ws.Range("I93").Value = 1680.5
$ws.Range("K93").Value = 1680.5
$ws.Range("M93").Value = -432.5
$ws.Range("H95").Value = 40344
$ws.Range("J95").Value = 40344
$ws.Range("L95").Value = 40344
$ws.Range("N95").Value = -45836
$ws.Range("H100").Value = 3939.2
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 2644
$ws.Range("I113").Value = 2710.9167
$ws.Range("J113").Value = 1038
$ws.Range("K113").Value = 2710.9167
$ws.Range("L113").Value = 1038
$ws.Range("M113").Value = -540.9167000000002
$ws.Range("N113").Value = -5378
$ws.Range("H122").Value = 10834.685
$ws.Range("I122").Value = 5019.4546
$ws.Range("J122").Value = 18830.625
$ws.Range("K122").Value = 15058.3638
$ws.Range("L122").Value = 56491.875
$ws.Range("M122").Value = -12608.3638
$ws.Range("N122").Value = -61391.875
$ws.Range("H126").Value = 71434570
$ws.Range("I126").Value = 90914296
$ws.Range("K126").Value = 272742888
$ws.Range("M126").Value = -272740418
$ws.Range("H132").Value = 5993.7754
$ws.Range("I132").Value = 5714.357
$ws.Range("K132").Value = 17143.071
$ws.Range("M132").Value = -14613.071
$ws.Range("H136").Value = 6738.35
$ws.Range("I136").Value = 8204.166999999999
$ws.Range("J136").Value = 4539.625
$ws.Range("K136").Value = 24612.501
$ws.Range("L136").Value = 13618.875
$ws.Range("M136").Value = -22062.501
$ws.Range("N136").Value = -18718.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3534
$ws.Range("I81").Value = 2834.8
$ws.Range("J81").Value = 10526
$ws.Range("K81").Value = 5669.6
$ws.Range("L81").Value = 21052
$ws.Range("M81").Value = -4608.6
$ws.Range("N81").Value = -23174
$ws.Range("H84").Value = 3534
$ws.Range("I84").Value = 2834.8
$ws.Range("J84").Value = 10526
$ws.Range("K84").Value = 28348
$ws.Range("L84").Value = 105260
$ws.Range("M84").Value = -23044
$ws.Range("N84").Value = -115868
$ws.Range("H96").Value = 1325322.4
$ws.Range("I96").Value = 3089916.2
$ws.Range("K96").Value = 3089916.2
$ws.Range("M96").Value = -3088543.2
$ws.Range("H107").Value = 1777.4286
$ws.Range("I107").Value = 1791.5
$ws.Range("J107").Value = 1758.6666
$ws.Range("K107").Value = 5374.5
$ws.Range("L107").Value = 5275.9998
$ws.Range("M107").Value = -3454.5
$ws.Range("N107").Value = -9115.9998
$ws.Range("H132").Value = 2490.2856
$ws.Range("I132").Value = 1766.0975
$ws.Range("J132").Value = 6201.75
$ws.Range("K132").Value = 5298.2925
$ws.Range("L132").Value = 18605.25
$ws.Range("M132").Value = -2768.2925
$ws.Range("N132").Value = -23665.25
$ws.Range("H136").Value = 5641.3477
$ws.Range("I136").Value = 4558.6294
$ws.Range("J136").Value = 7179.9473
$ws.Range("K136").Value = 13675.8882
$ws.Range("L136").Value = 21539.8419
$ws.Range("M136").Value = -11125.8882
$ws.Range("N136").Value = -26639.8419
